# Insert a new weekly price record for Coliflor (row 508) in the daily
# logic subset sheet. This pushes the existing rows 508..609 down to
# 509..610, matching the target diff (dimension grows from R609 to R610).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 508, shifting rows 508:609 down to 509:610.
$ws.Rows("508:508").Insert()

# Populate the newly inserted row 508 with the new record's data.
$ws.Range("A508").Value = 7
$ws.Range("B508").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C508").Value = "Ñuble"
$ws.Range("D508").Value = 45244
$ws.Range("D508").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E508").Value = 16
$ws.Range("F508").Value = 100112008
$ws.Range("G508").Value = "Coliflor"
$ws.Range("H508").Value = "Sin especificar"
$ws.Range("I508").Value = "Primera"
$ws.Range("J508").Value = 300
$ws.Range("K508").Value = 1300
$ws.Range("L508").Value = 1300
$ws.Range("M508").Value = 1300
$ws.Range("N508").Value = "$/unidad"
$ws.Range("O508").Value = "Región del Maule"
$ws.Range("P508").Value = 1300
$ws.Range("Q508").Value = 1
$ws.Range("R508").Value = "Hortaliza"
